# Update cryptocurrency price/volume figures (refresh pulled 2023-02-07).
# Values are entered with a leading apostrophe so they stay text cells
# (matching the workbook's existing inline-string / General-format cells)
# instead of being auto-converted to numbers/percentages by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.09"
$ws.Range("E2").Value = "'0.49%"
$ws.Range("D3").Value = "'44.16"
$ws.Range("E3").Value = "'-0.18%"
$ws.Range("D4").Value = "'5.557"
$ws.Range("E4").Value = "'-0.34%"
$ws.Range("D5").Value = "'0.08146"
$ws.Range("E5").Value = "'0.59%"
$ws.Range("D6").Value = "'2.059"
$ws.Range("E6").Value = "'4.05%"
$ws.Range("D7").Value = "'0.9747"
$ws.Range("E7").Value = "'2.42%"
$ws.Range("D8").Value = "'0.1104"
$ws.Range("E8").Value = "'-5.83%"
$ws.Range("D9").Value = "'0.1891"
$ws.Range("E9").Value = "'2.17%"
$ws.Range("D10").Value = "'10.13"
$ws.Range("E10").Value = "'-1.30%"
$ws.Range("D11").Value = "'0.09961"
$ws.Range("E11").Value = "'1.40%"
$ws.Range("D12").Value = "'0.04730"
$ws.Range("E12").Value = "'0.19%"
$ws.Range("D13").Value = "'0.1057"
$ws.Range("E13").Value = "'-1.18%"
$ws.Range("D14").Value = "'0.001274"
$ws.Range("E14").Value = "'-0.89%"
$ws.Range("D15").Value = "'0.04102"
$ws.Range("E15").Value = "'-3.07%"
$ws.Range("D16").Value = "'0.006096"
$ws.Range("E16").Value = "'3.02%"
$ws.Range("E17").Value = "'-0.87%"
$ws.Range("D18").Value = "'4.432"
$ws.Range("E18").Value = "'2.24%"
$ws.Range("E19").Value = "'1.65%"
$ws.Range("D20").Value = "'0.3349"
$ws.Range("E20").Value = "'-3.55%"
$ws.Range("E21").Value = "'-2.16%"
$ws.Range("D22").Value = "'0.2570"
$ws.Range("E22").Value = "'2.60%"
$ws.Range("E23").Value = "'3.92%"
$ws.Range("D24").Value = "'0.004377"
$ws.Range("E24").Value = "'0.82%"
$ws.Range("D25").Value = "'0.0001277"
$ws.Range("E25").Value = "'7.19%"
$ws.Range("D26").Value = "'0.0003736"
$ws.Range("E26").Value = "'-6.10%"
$ws.Range("D38").Value = "'0.02691"
$ws.Range("E38").Value = "'1.08%"
$ws.Range("D39").Value = "'0.05639"
$ws.Range("E39").Value = "'1.71%"
$ws.Range("D40").Value = "'0.007637"
$ws.Range("E40").Value = "'0.97%"
$ws.Range("D41").Value = "'0.1416"
$ws.Range("E41").Value = "'0.49%"
$ws.Range("E42").Value = "'-6.51%"
$ws.Range("D43").Value = "'0.001954"
$ws.Range("E43").Value = "'-3.14%"
$ws.Range("D44").Value = "'0.008307"
$ws.Range("E44").Value = "'-6.81%"
$ws.Range("E45").Value = "'-5.87%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.27%"
$ws.Range("D47").Value = "'0.0005794"
$ws.Range("E47").Value = "'-0.31%"
$ws.Range("E48").Value = "'10.71%"
$ws.Range("D49").Value = "'0.003534"
$ws.Range("E49").Value = "'0.67%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.27%"
